# Auto-generated PowerShell COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style taken from an untouched data cell (B2) which carries the
# workbook's default/no-op style (no explicit s= attribute). We reapply it
# after forcing NumberFormat to Text so numeric-looking strings (e.g. '5.19')
# are preserved verbatim as text, matching the source inline strings, without
# leaving a stray custom number-format style behind on the cell.
$refStyle = $ws.Range("B2").Style

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $refStyle
}

# Row 2
Set-TextCell "D2" '58.804.88'
Set-TextCell "E2" '  +1.05%  '

# Row 3
Set-TextCell "D3" '2.501.99'
Set-TextCell "E3" '  +0.81%  '

# Row 4
Set-TextCell "E4" '  +0.09%  '

# Row 5
Set-TextCell "D5" '532.88'
Set-TextCell "E5" '  +3.18%  '

# Row 6
Set-TextCell "D6" '133.70'
Set-TextCell "E6" '  +1.55%  '

# Row 7
Set-TextCell "E7" '  +0.42%  '

# Row 8
Set-TextCell "D8" '0.570'
Set-TextCell "E8" '  +3.06%  '

# Row 9
Set-TextCell "D9" '2.503.93'
Set-TextCell "E9" '  -0.47%  '

# Row 10
Set-TextCell "E10" '  +1.99%  '

# Row 11
Set-TextCell "E11" '  -2.59%  '

# Row 12
Set-TextCell "D12" '5.19'
Set-TextCell "E12" '  -0.48%  '

# Row 13
Set-TextCell "E13" '  -0.78%  '

# Row 14
Set-TextCell "D14" '2.944.80'
Set-TextCell "E14" '  +1.10%  '

# Row 15
Set-TextCell "D15" '58.664.11'
Set-TextCell "E15" '  +1.00%  '

# Row 16
Set-TextCell "D16" '22.26'
Set-TextCell "E16" '  +0.58%  '

# Row 17
Set-TextCell "E17" '  +0.91%  '

# Row 18
Set-TextCell "D18" '2.494.67'
Set-TextCell "E18" '  +0.70%  '

# Row 19
Set-TextCell "D19" '10.58'
Set-TextCell "E19" '  -0.65%  '

# Row 20
Set-TextCell "E20" '  +2.18%  '

# Row 21
Set-TextCell "D21" '320.79'
Set-TextCell "E21" '  +0.16%  '

# Row 22
Set-TextCell "D22" '6.20'
Set-TextCell "E22" '  +3.51%  '

# Row 23
Set-TextCell "D23" '1.00'
Set-TextCell "E23" '  -0.04%  '

# Row 24
Set-TextCell "E24" '  +4.54%  '

# Row 25
Set-TextCell "E25" '  +0.85%  '

# Row 26
Set-TextCell "E26" '  +0.13%  '

# Row 27
Set-TextCell "E27" '  -1.00%  '

# Row 28
Set-TextCell "E28" '  +1.62%  '

# Row 29
Set-TextCell "D29" '172.70'
Set-TextCell "E29" '  +1.82%  '

# Row 30
Set-TextCell "D30" '0.0₃0755'
Set-TextCell "E30" '  +1.79%  '

# Row 31
Set-TextCell "E31" '  +2.50%  '

# Row 32
Set-TextCell "D32" '6.27'
Set-TextCell "E32" '  -0.16%  '

# Row 33
Set-TextCell "E33" '  -1.06%  '

# Row 35
Set-TextCell "E35" '  +0.60%  '

# Row 36
Set-TextCell "E36" '  +0.44%  '

# Row 37
Set-TextCell "E37" '  -4.49%  '

# Row 38
Set-TextCell "D38" '3.94'
Set-TextCell "E38" '  -0.32%  '

# Row 39
Set-TextCell "B39" 'Stacks'
Set-TextCell "C39" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D39" '1.51'
Set-TextCell "E39" '  +2.71%  '

# Row 40
Set-TextCell "B40" 'SuiNetwork'
Set-TextCell "C40" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell "D40" '0.832'
Set-TextCell "E40" '  +6.67%  '

# Row 41
Set-TextCell "D41" '36.37'
Set-TextCell "E41" '  -0.91%  '

# Row 42
Set-TextCell "D42" '3.46'
Set-TextCell "E42" '  +1.78%  '

# Row 43
Set-TextCell "D43" '274.72'
Set-TextCell "E43" '  -0.15%  '

# Row 44
Set-TextCell "E44" '  +6.80%  '

# Row 45
Set-TextCell "E45" '  -1.96%  '

# Row 46
Set-TextCell "D46" '0.591'
Set-TextCell "E46" '  -0.50%  '

# Row 47
Set-TextCell "D47" '0.0932'
Set-TextCell "E47" '  +1.74%  '

# Row 48
Set-TextCell "D48" '0.0508'
Set-TextCell "E48" '  +3.06%  '

# Row 49
Set-TextCell "E49" '  +2.09%  '

# Row 50
Set-TextCell "D50" '16.75'
Set-TextCell "E50" '  -0.55%  '

# Row 51
Set-TextCell "D51" '1.753.71'
Set-TextCell "E51" '  +1.43%  '
